$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (shifts old C -> D, keeps its formatting)
$ws.Columns("C:C").Insert()

# Fill the new "Turha"/"turhuus" column (C), matching the plain header/body font (style 1)
$ws.Range("C1:C6").Font.Size = 15
$ws.Range("C1").Value = "Turha"
$ws.Range("C2").Value = "turhuus"
$ws.Range("C3").Value = "turhuus"
$ws.Range("C4").Value = "turhuus"
$ws.Range("C5").Value = "turhuus"
$ws.Range("C6").Value = "turhuus"

# Fill the new "Uusi "/"Testaa" column (E), same plain font (style 1)
$ws.Range("E1:E6").Font.Size = 15
$ws.Range("E1").Value = "Uusi "
$ws.Range("E2").Value = "Testaa"
$ws.Range("E3").Value = "Testaa"
$ws.Range("E4").Value = "Testaa"
$ws.Range("E5").Value = "Testaa"
$ws.Range("E6").Value = "Testaa"

# Widen the new column C (~18.27 characters, as in the target workbook)
$ws.Columns("C:C").ColumnWidth = 17.5

# Update the window/selection state to match
$ws.Range("C10").Select()
